$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that get the "On Order" status (Neutral cell style) in new column G
$onOrderRows = @(6, 7, 9, 16, 17, 19, 25, 32)
foreach ($r in $onOrderRows) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = "On Order"
    $cell.Style = "Neutral"
}

# Rows that get the "Need" status (Accent1 cell style) in new column G
$needRows = @(8, 12, 15, 30)
foreach ($r in $needRows) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = "Need"
    $cell.Style = "Accent1"
}

# Restore the selected cell as left by the author
$ws.Range("J20").Select()
